# "Failure to Launch" sheet: the "Launch Failure Date" column (A2:A11) was
# stored as real date serials formatted with a custom "M/d/yyyy" number
# format. Convert those dates to plain text (keeping the same displayed
# value) so the column no longer depends on the custom date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Failure to Launch")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# Capture the currently-displayed date text for every row BEFORE touching
# the number format, since changing the format would change what .Text
# reports.
$dateText = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null) {
        $dateText[$r] = $cell.Text
    }
}

# Switch the column to plain text and write the captured date strings back
# in, so the cells now hold literal text like "10/23/2035" instead of a
# numeric date serial.
$ws.Range("A2:A$lastRow").NumberFormat = "@"
foreach ($r in $dateText.Keys) {
    $ws.Cells.Item($r, 1).Value = $dateText[$r]
}
